$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.322.80'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.882.40'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.60'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.44%  '
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.36'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.354'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '53.41'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0740'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0970'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.44'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.24%  '
$ws.Range('D14').Value = '2.159.00'
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.761'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.89'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.25%  '
$ws.Range('D17').Value = '1.877.64'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('D18').Value = '35.483.32'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.85'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').Value = '0.0₃0820'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '243.50'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.76'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.94'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('E24').Value = '  +8.76%  '
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.15'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -5.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.08'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.52'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.30'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.23%  '
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.73'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +10.14%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.99'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.77%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.25'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0586'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.13'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.65%  '
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.94'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0718'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +10.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.42'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0217'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.76'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -7.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.07'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.45%  '
$ws.Range('D45').Value = '1.303.36'
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('E46').Value = '  -2.95%  '
$ws.Range('E47').Value = '  +6.87%  '
$ws.Range('E48').Value = '  -1.96%  '
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.16'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.23'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -5.49%  '
